$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet currently has data through row 385 (A1:D385).
# Extend the table through row 464, matching the existing formatting
# (column A uses style index 2 - bold border, centered, custom date format)
# by copying the formatting of the last existing row down across the new rows.
$ws.Range("A385:D385").Copy() | Out-Null
$ws.Range("A386:D464").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New daily data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$data = @(
    @(386, 44460, 0, 1, 62.34413965087282),
    @(387, 44461, 0, 1, 62.34413965087282),
    @(388, 44462, 0, 1, 62.34413965087282),
    @(389, 44463, 0, 1, 62.34413965087282),
    @(390, 44464, 0, 0, 0),
    @(391, 44465, 0, 0, 0),
    @(392, 44466, 0, 0, 0),
    @(393, 44467, 1, 1, 62.34413965087282),
    @(394, 44468, 0, 1, 62.34413965087282),
    @(395, 44469, 0, 1, 62.34413965087282),
    @(396, 44470, 0, 1, 62.34413965087282),
    @(397, 44471, 0, 1, 62.34413965087282),
    @(398, 44472, 0, 1, 62.34413965087282),
    @(399, 44473, 0, 1, 62.34413965087282),
    @(400, 44474, 0, 0, 0),
    @(401, 44475, 1, 1, 62.34413965087282),
    @(402, 44476, 0, 1, 62.34413965087282),
    @(403, 44477, 1, 2, 124.6882793017456),
    @(404, 44478, 0, 2, 124.6882793017456),
    @(405, 44479, 0, 2, 124.6882793017456),
    @(406, 44480, 0, 2, 124.6882793017456),
    @(407, 44481, 0, 2, 124.6882793017456),
    @(408, 44482, 0, 1, 62.34413965087282),
    @(409, 44483, 0, 1, 62.34413965087282),
    @(410, 44484, 0, 0, 0),
    @(411, 44485, 0, 0, 0),
    @(412, 44486, 0, 0, 0),
    @(413, 44487, 0, 0, 0),
    @(414, 44488, 0, 0, 0),
    @(415, 44489, 0, 0, 0),
    @(416, 44490, 0, 0, 0),
    @(417, 44491, 0, 0, 0),
    @(418, 44492, 0, 0, 0),
    @(419, 44493, 0, 0, 0),
    @(420, 44494, 0, 0, 0),
    @(421, 44495, 0, 0, 0),
    @(422, 44496, 0, 0, 0),
    @(423, 44497, 0, 0, 0),
    @(424, 44498, 0, 0, 0),
    @(425, 44499, 0, 0, 0),
    @(426, 44500, 0, 0, 0),
    @(427, 44501, 0, 0, 0),
    @(428, 44502, 0, 0, 0),
    @(429, 44503, 0, 0, 0),
    @(430, 44504, 0, 0, 0),
    @(431, 44505, 0, 0, 0),
    @(432, 44506, 0, 0, 0),
    @(433, 44507, 0, 0, 0),
    @(434, 44508, 0, 0, 0),
    @(435, 44509, 0, 0, 0),
    @(436, 44510, 0, 0, 0),
    @(437, 44511, 0, 0, 0),
    @(438, 44512, 0, 0, 0),
    @(439, 44513, 0, 0, 0),
    @(440, 44514, 0, 0, 0),
    @(441, 44515, 0, 0, 0),
    @(442, 44516, 0, 0, 0),
    @(443, 44517, 0, 0, 0),
    @(444, 44518, 0, 0, 0),
    @(445, 44519, 0, 0, 0),
    @(446, 44520, 0, 0, 0),
    @(447, 44521, 0, 0, 0),
    @(448, 44522, 0, 0, 0),
    @(449, 44523, 0, 0, 0),
    @(450, 44524, 1, 1, 62.34413965087282),
    @(451, 44525, 0, 1, 62.34413965087282),
    @(452, 44526, 0, 1, 62.34413965087282),
    @(453, 44527, 0, 1, 62.34413965087282),
    @(454, 44528, 0, 1, 62.34413965087282),
    @(455, 44529, 0, 1, 62.34413965087282),
    @(456, 44530, 0, 1, 62.34413965087282),
    @(457, 44531, 0, 0, 0),
    @(458, 44532, 0, 0, 0),
    @(459, 44533, 0, 0, 0),
    @(460, 44534, 0, 0, 0),
    @(461, 44535, 0, 0, 0),
    @(462, 44536, 0, 0, 0),
    @(463, 44537, 0, 0, 0),
    @(464, 44538, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
